# Re-sorts each "Successor N" block (5 country rows) alphabetically by
# country name (column A) and refreshes each block's "Big-U" value (E col
# on the block header row), matching the recomputed alloy_transform pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Successors (Test Results)")

$blocks = @(
    @{ HeaderRow=1; BigU=12.37110849056604; Rows=@(@{Row=3; A="Atlantis"; Vals=@(100,698,2000,1,0,0,1,0,0)}, @{Row=4; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=5; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=6; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=7; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=9; BigU=12.29825187969925; Rows=@(@{Row=11; A="Atlantis"; Vals=@(100,696,2000,2,0,0,2,0,0)}, @{Row=12; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=13; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=14; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=15; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=17; BigU=12.2259691011236; Rows=@(@{Row=19; A="Atlantis"; Vals=@(100,694,2000,3,0,0,3,0,0)}, @{Row=20; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=21; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=22; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=23; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=25; BigU=12.15425373134329; Rows=@(@{Row=27; A="Atlantis"; Vals=@(100,692,2000,4,0,0,4,0,0)}, @{Row=28; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=29; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=30; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=31; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=33; BigU=12.08309944237918; Rows=@(@{Row=35; A="Atlantis"; Vals=@(100,690,2000,5,0,0,5,0,0)}, @{Row=36; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=37; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=38; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=39; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=41; BigU=12.0125; Rows=@(@{Row=43; A="Atlantis"; Vals=@(100,688,2000,6,0,0,6,0,0)}, @{Row=44; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=45; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=46; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=47; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=49; BigU=11.94244926199263; Rows=@(@{Row=51; A="Atlantis"; Vals=@(100,686,2000,7,0,0,7,0,0)}, @{Row=52; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=53; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=54; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=55; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=57; BigU=11.87294117647059; Rows=@(@{Row=59; A="Atlantis"; Vals=@(100,684,2000,8,0,0,8,0,0)}, @{Row=60; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=61; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=62; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=63; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=65; BigU=11.80396978021979; Rows=@(@{Row=67; A="Atlantis"; Vals=@(100,682,2000,9,0,0,9,0,0)}, @{Row=68; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=69; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=70; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=71; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) },
    @{ HeaderRow=73; BigU=11.73552919708029; Rows=@(@{Row=75; A="Atlantis"; Vals=@(100,680,2000,10,0,0,10,0,0)}, @{Row=76; A="Brobdingnag"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=77; A="Carpania"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=78; A="Dinotopia"; Vals=@(1,0,0,0,0,0,0,0,0)}, @{Row=79; A="Erewhon"; Vals=@(1,0,0,0,0,0,0,0,0)}) }
)

foreach ($block in $blocks) {
    $ws.Cells.Item($block.HeaderRow, 5).Value = $block.BigU
    foreach ($row in $block.Rows) {
        $ws.Cells.Item($row.Row, 1).Value = $row.A
        for ($col = 2; $col -le 10; $col++) {
            $ws.Cells.Item($row.Row, $col).Value = $row.Vals[$col - 2]
        }
    }
}